$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 212, shifting existing rows
# 212-278 down to 214-280 (Excel preserves their values/styles automatically).
$ws.Rows("212:213").Insert()

# Populate new row 212 with the new price record.
$ws.Cells.Item(212, 1).Value = 8
$ws.Cells.Item(212, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(212, 3).Value = "Coquimbo"
$ws.Cells.Item(212, 4).Value = 44809
$ws.Cells.Item(212, 5).Value = 4
$ws.Cells.Item(212, 6).Value = 100112021
$ws.Cells.Item(212, 7).Value = "Ají"
$ws.Cells.Item(212, 8).Value = "Inferno"
$ws.Cells.Item(212, 9).Value = "Primera"
$ws.Cells.Item(212, 10).Value = 520
$ws.Cells.Item(212, 11).Value = 16000
$ws.Cells.Item(212, 12).Value = 17000
$ws.Cells.Item(212, 13).Value = 16500
$ws.Cells.Item(212, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(212, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(212, 16).Value = 1650
$ws.Cells.Item(212, 17).Value = 10
$ws.Cells.Item(212, 18).Value = "Hortaliza"

# Populate new row 213 with the new price record.
$ws.Cells.Item(213, 1).Value = 8
$ws.Cells.Item(213, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(213, 3).Value = "Coquimbo"
$ws.Cells.Item(213, 4).Value = 44809
$ws.Cells.Item(213, 5).Value = 4
$ws.Cells.Item(213, 6).Value = 100112021
$ws.Cells.Item(213, 7).Value = "Ají"
$ws.Cells.Item(213, 8).Value = "Inferno"
$ws.Cells.Item(213, 9).Value = "Segunda"
$ws.Cells.Item(213, 10).Value = 400
$ws.Cells.Item(213, 11).Value = 10000
$ws.Cells.Item(213, 12).Value = 11000
$ws.Cells.Item(213, 13).Value = 10500
$ws.Cells.Item(213, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(213, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(213, 16).Value = 1050
$ws.Cells.Item(213, 17).Value = 10
$ws.Cells.Item(213, 18).Value = "Hortaliza"
